$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: replace month text with numeric month index
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 3

# Column C: shift studentId values down and fill new rows
$ws.Range("C2").Value = "st45"
$ws.Range("C3").Value = "st55"

# New rows 4-6: paginate additional student marks records
$ws.Range("A4").Value = 45
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "st65"

$ws.Range("A5").Value = 69
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "st75"

$ws.Range("A6").Value = 56
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "st85"

$ws.Range("C6").Select()
